# Cronograma.xlsx update script
# - Inserts a new task row ("Regresion logistica") after row 10 (old row 11
#   "Implementacion Q-Learning" and everything below through row 21 shifts
#   down by one row).
# - Adds two new date columns (AG, AH) for 2020-07-14 with new hour entries.
# - Updates the running total text in A1 and the raw hour-sum string in (the
#   row that was) A24.
# - Moves the selection back to a non-scrolled, different active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 11 - this shifts rows 11-23 down by one
#    (rows 11-20 carried data, 21-23 were blank, so the formerly-row-24
#    "total hours" line lands on row 25 and needs to be moved back up).
# ---------------------------------------------------------------------------
[void]$ws.Rows("11:11").Insert()

# Clean up the formatting Excel auto-copied down into the new blank row 11
# from row 10 (only column A and AH should carry content/format here).
[void]$ws.Range("L11:N11").Clear()

# New row 11 content: "Regresion logistica" task with a single "2 h." entry
# in the newly added AH column. Column A already inherited the correct
# "task name" style (s=3) from row 10 via the row insert above, so only the
# value needs to be set there. AH11 is a brand new column, so its format
# (matching the other "N h." entry cells) is copied in before the value.
$ws.Range("A11").Value = "Regresión logística"

[void]$ws.Range("L10").Copy()
[void]$ws.Range("AH11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AH11").Value = "2 h."

# ---------------------------------------------------------------------------
# 2. Move the "total hours" row (now on row 25 after the insert) back to
#    row 24, matching the target layout, and update its text.
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = $ws.Range("A25").Value
[void]$ws.Rows("25:25").Delete()
$ws.Range("A24").Value = "4+3+1.5+4+4+3+2+3.5+4.5+4.5+5+1+5+3+4+5+4+4+3+4+3.5+8+3.5+5+5+5+4+3+3.5+4+2.5+5.5"

# ---------------------------------------------------------------------------
# 3. Update the header total-hours label in A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Total horas: 124.5"

# ---------------------------------------------------------------------------
# 4. Add the two new date columns (AG and AH) with the 2020-07-14 date in
#    row 2, matching the style used by the existing date cells.
# ---------------------------------------------------------------------------
[void]$ws.Range("AF2").Copy()
[void]$ws.Range("AG2:AH2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AG2").Value = 44026
$ws.Range("AH2").Value = 44026

# ---------------------------------------------------------------------------
# 5. Populate the new hour entries that fall in the new AG/AH columns and
#    the newly freed-up cells after the row shift.
# ---------------------------------------------------------------------------

# Row 5 ("Generacion de graficas") -> AH5 = "1 h."
[void]$ws.Range("C5").Copy()
[void]$ws.Range("AH5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AH5").Value = "1 h."

# Row 13 ("Estructuras basicas (estados y agente)") -> AH13 blank placeholder
# cell with a right-aligned, unfilled style (new cellXfs entry).
$ws.Range("AH13").HorizontalAlignment = -4152

# Row 14 ("Algoritmo Q-learning") -> AG14 = "2.5 h."
[void]$ws.Range("AF14").Copy()
[void]$ws.Range("AG14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AG14").Value = "2.5 h."

# Row 18 ("Servidor") -> AH18 = "2.5 h."
[void]$ws.Range("AC18").Copy()
[void]$ws.Range("AH18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AH18").Value = "2.5 h."

# ---------------------------------------------------------------------------
# 6. Update the view: select L31 (also clears any scrolled top-left cell).
# ---------------------------------------------------------------------------
[void]$ws.Range("L31").Select()

Write-Host "Edit script completed"
